$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 112.5
$ws.Range("J2").Value = 150
$ws.Range("L2").Value = 150
$ws.Range("N2").Value = -376

$ws.Range("H18").Value = 211.81818
$ws.Range("I18").Value = 163
$ws.Range("K18").Value = 163
$ws.Range("M18").Value = 121

$ws.Range("H40").Value = 1496.0834
$ws.Range("I40").Value = 1071.8334
$ws.Range("J40").Value = 1920.3334
$ws.Range("K40").Value = 1071.8334
$ws.Range("L40").Value = 1920.3334
$ws.Range("M40").Value = -896.8334
$ws.Range("N40").Value = -2270.3334

$ws.Range("H55").Value = 430
$ws.Range("I55").Value = 668.75
$ws.Range("J55").Value = 157.14285
$ws.Range("K55").Value = 668.75
$ws.Range("L55").Value = 157.14285
$ws.Range("M55").Value = -454.75
$ws.Range("N55").Value = -585.14285

$ws.Range("H70").Value = 881.5714
$ws.Range("I70").Value = 714
$ws.Range("K70").Value = 2142
$ws.Range("M70").Value = -1872

$ws.Range("H73").Value = 881.5714
$ws.Range("I73").Value = 714
$ws.Range("K73").Value = 2142
$ws.Range("M73").Value = -1206

$ws.Range("H136").Value = 23750
$ws.Range("J136").Value = 23750
$ws.Range("L136").Value = 23750
$ws.Range("N136").Value = -33950

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2238.0476
$ws.Range("I61").Value = 2014.5625
$ws.Range("J61").Value = 2953.2
$ws.Range("K61").Value = 2014.5625
$ws.Range("L61").Value = 2953.2
$ws.Range("M61").Value = -1802.5625
$ws.Range("N61").Value = -3377.2

$ws.Range("H132").Value = 17042.576
$ws.Range("I132").Value = 1738.1724
$ws.Range("K132").Value = 5214.5172
$ws.Range("M132").Value = -2684.5172

$ws.Range("H136").Value = 2238.0476
$ws.Range("I136").Value = 2014.5625
$ws.Range("J136").Value = 2953.2
$ws.Range("K136").Value = 6043.6875
$ws.Range("L136").Value = 8859.599999999999
$ws.Range("M136").Value = -3493.6875
$ws.Range("N136").Value = -13959.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4583.2856
$ws.Range("I20").Value = 5998.6
$ws.Range("K20").Value = 5998.6
$ws.Range("M20").Value = -5751.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 7999
$ws.Range("I14").Value = 7999
$ws.Range("K14").Value = 7999
$ws.Range("M14").Value = -7829

$ws.Range("H21").Value = 14625
$ws.Range("I21").Value = 12000
$ws.Range("J21").Value = 15000
$ws.Range("K21").Value = 12000
$ws.Range("L21").Value = 15000
$ws.Range("M21").Value = -11765
$ws.Range("N21").Value = -15470

$ws.Range("H58").Value = 9140.271000000001
$ws.Range("I58").Value = 631.093
$ws.Range("J58").Value = 32008.688
$ws.Range("K58").Value = 631.093
$ws.Range("L58").Value = 32008.688
$ws.Range("M58").Value = -428.093
$ws.Range("N58").Value = -32414.688

$ws.Range("H59").Value = 26533.334
$ws.Range("J59").Value = 26533.334
$ws.Range("L59").Value = 26533.334
$ws.Range("N59").Value = -28823.334

$ws.Range("H132").Value = 14369.452
$ws.Range("I132").Value = 19167.725
$ws.Range("J132").Value = 3665.6155
$ws.Range("K132").Value = 57503.175
$ws.Range("L132").Value = 10996.8465
$ws.Range("M132").Value = -54973.175
$ws.Range("N132").Value = -16056.8465

$ws.Range("H134").Value = 652.6429000000001
$ws.Range("I134").Value = 574.8
$ws.Range("K134").Value = 1724.4
$ws.Range("M134").Value = 810.6000000000001

$ws.Range("H136").Value = 9140.271000000001
$ws.Range("I136").Value = 631.093
$ws.Range("J136").Value = 32008.688
$ws.Range("K136").Value = 1893.279
$ws.Range("L136").Value = 96026.064
$ws.Range("M136").Value = 656.721
$ws.Range("N136").Value = -101126.064

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5812
$ws.Range("I3").Value = 1790
$ws.Range("K3").Value = 5370
$ws.Range("M3").Value = -5258

$ws.Range("H107").Value = 4268.4
$ws.Range("J107").Value = 280.5
$ws.Range("L107").Value = 841.5
$ws.Range("N107").Value = -4681.5

$ws.Range("H131").Value = 766.02
$ws.Range("J131").Value = 776.04083
$ws.Range("L131").Value = 2328.12249
$ws.Range("N131").Value = -12408.12249

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 38462852
$ws.Range("I102").Value = 41667980
$ws.Range("J102").Value = 1314
$ws.Range("K102").Value = 41667980
$ws.Range("L102").Value = 1314
$ws.Range("M102").Value = -41666358
$ws.Range("N102").Value = -4558

$ws.Range("H114").Value = 43760
$ws.Range("J114").Value = 43760
$ws.Range("L114").Value = 43760
$ws.Range("N114").Value = -52438

$ws.Range("H132").Value = 22797
$ws.Range("I132").Value = 4676.05
$ws.Range("J132").Value = 74571.14
$ws.Range("K132").Value = 14028.15
$ws.Range("L132").Value = 223713.42
$ws.Range("M132").Value = -11498.15
$ws.Range("N132").Value = -228773.42

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 864.41174
$ws.Range("I46").Value = 999.1667
$ws.Range("J46").Value = 790.9091
$ws.Range("K46").Value = 999.1667
$ws.Range("L46").Value = 790.9091
$ws.Range("M46").Value = -811.1667
$ws.Range("N46").Value = -1166.9091

$ws.Range("H68").Value = 1977.375
$ws.Range("I68").Value = 1313.9
$ws.Range("J68").Value = 3083.1667
$ws.Range("K68").Value = 1313.9
$ws.Range("L68").Value = 3083.1667
$ws.Range("M68").Value = -564.9000000000001
$ws.Range("N68").Value = -4581.1667

$ws.Range("H71").Value = 1977.375
$ws.Range("I71").Value = 1313.9
$ws.Range("J71").Value = 3083.1667
$ws.Range("K71").Value = 6569.5
$ws.Range("L71").Value = 15415.8335
$ws.Range("M71").Value = -2825.5
$ws.Range("N71").Value = -22903.8335

$ws.Range("H93").Value = 2743.125
$ws.Range("I93").Value = 2590.5
$ws.Range("J93").Value = 3201
$ws.Range("K93").Value = 2590.5
$ws.Range("L93").Value = 3201
$ws.Range("M93").Value = -1342.5
$ws.Range("N93").Value = -5697

$ws.Range("H136").Value = 18713.893
$ws.Range("I136").Value = 23463.137
$ws.Range("K136").Value = 70389.41099999999
$ws.Range("M136").Value = -67839.41099999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 2853.3333
$ws.Range("J22").Value = 2853.3333
$ws.Range("L22").Value = 2853.3333
$ws.Range("N22").Value = -3439.3333

$ws.Range("H54").Value = 17500
$ws.Range("J54").Value = 17500
$ws.Range("L54").Value = 17500
$ws.Range("N54").Value = -18540

$ws.Range("H81").Value = 100001670
$ws.Range("I81").Value = 2104.5715
$ws.Range("J81").Value = 333334000
$ws.Range("K81").Value = 4209.143
$ws.Range("L81").Value = 666668000
$ws.Range("M81").Value = -3148.143
$ws.Range("N81").Value = -666670122

$ws.Range("H84").Value = 100001670
$ws.Range("I84").Value = 2104.5715
$ws.Range("J84").Value = 333334000
$ws.Range("K84").Value = 21045.715
$ws.Range("L84").Value = 3333340000
$ws.Range("M84").Value = -15741.715
$ws.Range("N84").Value = -3333350608

$ws.Range("H132").Value = 1009.775
$ws.Range("I132").Value = 709.1613
$ws.Range("J132").Value = 2045.2222
$ws.Range("K132").Value = 2127.4839
$ws.Range("L132").Value = 6135.6666
$ws.Range("M132").Value = 402.5160999999998
$ws.Range("N132").Value = -11195.6666
